$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 4 (pushes old rows 4..25 down to 5..26),
# carrying the existing formatting/validation/conditional-formatting
# ranges along with it (matches "Fixed new tab wait problem": a new
# "Register test" / registertest scenario row was added right after the
# "003" row).
$ws.Rows.Item(4).Insert()

# Populate the new row 4 with the "Register test" scenario.
# Leading apostrophe keeps the numeric-looking ID stored as text (like
# the other ID cells in the sheet, e.g. "001", "002", "003"...).
$ws.Cells.Item(4, 1).Value = "'003"
$ws.Cells.Item(4, 2).Value = "Register test"
$ws.Cells.Item(4, 3).Value = "ExpandTests"
$ws.Cells.Item(4, 4).Value = "registertest"
$ws.Cells.Item(4, 5).Value = $true

# Move the active selection to the new row's Run cell.
$ws.Range("E4").Select()
